# Trader Oracle Testing.xlsx - add Reports section, rework Login/Logout rows
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Drop everything below the Login block (rows 7-21) ---------------------
$ws.Range("A7:I21").EntireRow.Delete()

# --- Login section tweaks ---------------------------------------------------
# New test case replacing "not activated email" -> "not registered email"
$ws.Range("B5").Value = "Trader login with not registered email."
$ws.Range("C5").Value = "Error message is displayed."

# "valid credentials" row moves down one slot, Steps/Expected swap order
$ws.Range("B6").Value = "Trader login with valid credentials."
$ws.Range("C6").Value = "Login success and redirected to dashboard."

# --- New "Reports" section (rows 8-13) --------------------------------------
$ws.Range("A8").Value = "Reports"
$ws.Range("A8").Font.Bold = $true

$reportRows = @(
    @{ Row = 8;  Steps = "View Monthly Order Report";        Expected = "Monthly Order Report is displayed." },
    @{ Row = 9;  Steps = "View Products Stock Report";       Expected = "Products Stock Report is displayed." },
    @{ Row = 10; Steps = "View Payment Report";               Expected = "Payment Report is displayed." },
    @{ Row = 11; Steps = "View Weekly Order Report";          Expected = "Weekly Order Report is displayed." },
    @{ Row = 12; Steps = "View Daily Order Report";           Expected = "Daily Order Report is displayed." },
    @{ Row = 13; Steps = "View Review and Ratings Report";    Expected = "Review and Ratings report is displayed." }
)

foreach ($r in $reportRows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.Steps
    $ws.Cells.Item($r.Row, 3).Value = $r.Expected
    $ws.Cells.Item($r.Row, 4).Value = "Pass"
    $ws.Cells.Item($r.Row, 4).Font.Bold = $true
    $ws.Cells.Item($r.Row, 4).Font.Color = 5287936
}

# row13's section-label cell (column A) stays blank but keeps the bold style
$ws.Range("A13").Font.Bold = $true

# --- Logout section moves from row 21 to row 15 -----------------------------
$ws.Range("A15").Value = "Logout"
$ws.Range("A15").Font.Bold = $true

$ws.Range("B15").Value = "Trader Logout"
$ws.Range("C15").Value = "Trader gets logged out."
$ws.Range("D15").Value = "Pass"
$ws.Range("D15").Font.Bold = $true
$ws.Range("D15").Font.Color = 5287936

# --- Sheet view: selection moved back to B5, no frozen top-left cell -------
$ws.Range("B5").Select()
